# Heap Implementation containing all the methods of HEAP ADT like
# percolateUp, percolateDown -- add new rows/columns documenting three
# more problems (kth smallest in min-heap, merge K sorted arrays,
# median of an infinite stream) plus Code/Algo, T.C and S.C columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the two existing row-1 headers and the row-2 class name before
# they get shifted from D/E to G/H by the newly inserted columns.
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2
$oldE2 = $ws.Range("E2").Value2

# ---- Row 1 (headers) ----
$ws.Range("D1").Value = "Code/Algo"
$ws.Range("E1").Value = "T.C"
$ws.Range("F1").Value = "S.C"
$ws.Range("G1").Value = $oldD1
$ws.Range("H1").Value = $oldE1

# ---- Row 2 (Heap / Heap ADT) ----
$ws.Range("D2").Value = "c"
$ws.Range("E2").ClearContents()
$ws.Range("H2").Value = $oldE2

# ---- Row 3 (Finding kth Smallest element in Min Heap) ----
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Finding kth Smallest element in Min Heap"
$ws.Range("D3").Value = "a"
$ws.Range("E3").Value = "O(klogn), O(klogk)"
$ws.Range("F3").Value = "O(1),O(k)"

# ---- Row 4 (Merge K sorted arrays of size n each) ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Merge K sorted arrays of size n each"
$ws.Range("D4").Value = "a"
$ws.Range("E4").Value = "O(nk2),O(nklogk)"

# ---- Row 5 (Median of numbers in infinite stream of integers) ----
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Median of numbers in infinite stream of integers"
$ws.Range("D5").Value = "c"

# Resize the columns to fit their new contents.
$ws.Range("A1:H5").EntireColumn.AutoFit() | Out-Null

$ws.Range("D5").Select()
